$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 31; Time = "2023-12-06 17:14:33"; Cost = 0.0002 },
    @{ Row = 32; Time = "2023-12-06 17:15:42"; Cost = 0.004000000000000001 },
    @{ Row = 33; Time = "2023-12-06 17:16:11"; Cost = 0.0024 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Time
    $ws.Cells.Item($r.Row, 2).Value = $r.Cost
}
